# Updates the lattice multiplication practice table: each of the 15
# cells gets a new "AA x BB" problem, new partial-product digits, and
# new lattice-column leader digits. Cell (row, col) order matches the
# table's natural reading order.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11   # manual line break (w:br), matches Shift+Enter in Word

$t.Cell(1, 1).Range.Text = "53 x 85" + $nl + "  8    5" + $nl + "  ----" + $nl + "5|    |" + $nl + "3|    |"
$t.Cell(1, 2).Range.Text = "21 x 86" + $nl + "  8    6" + $nl + "  ----" + $nl + "2|    |" + $nl + "1|    |"
$t.Cell(1, 3).Range.Text = "16 x 89" + $nl + "  8    9" + $nl + "  ----" + $nl + "1|    |" + $nl + "6|    |"
$t.Cell(2, 1).Range.Text = "37 x 68" + $nl + "  6    8" + $nl + "  ----" + $nl + "3|    |" + $nl + "7|    |"
$t.Cell(2, 2).Range.Text = "71 x 66" + $nl + "  6    6" + $nl + "  ----" + $nl + "7|    |" + $nl + "1|    |"
$t.Cell(2, 3).Range.Text = "42 x 42" + $nl + "  4    2" + $nl + "  ----" + $nl + "4|    |" + $nl + "2|    |"
$t.Cell(3, 1).Range.Text = "65 x 68" + $nl + "  6    8" + $nl + "  ----" + $nl + "6|    |" + $nl + "5|    |"
$t.Cell(3, 2).Range.Text = "64 x 32" + $nl + "  3    2" + $nl + "  ----" + $nl + "6|    |" + $nl + "4|    |"
$t.Cell(3, 3).Range.Text = "86 x 65" + $nl + "  6    5" + $nl + "  ----" + $nl + "8|    |" + $nl + "6|    |"
$t.Cell(4, 1).Range.Text = "90 x 13" + $nl + "  1    3" + $nl + "  ----" + $nl + "9|    |" + $nl + "0|    |"
$t.Cell(4, 2).Range.Text = "60 x 32" + $nl + "  3    2" + $nl + "  ----" + $nl + "6|    |" + $nl + "0|    |"
$t.Cell(4, 3).Range.Text = "47 x 52" + $nl + "  5    2" + $nl + "  ----" + $nl + "4|    |" + $nl + "7|    |"
$t.Cell(5, 1).Range.Text = "89 x 30" + $nl + "  3    0" + $nl + "  ----" + $nl + "8|    |" + $nl + "9|    |"
$t.Cell(5, 2).Range.Text = "40 x 90" + $nl + "  9    0" + $nl + "  ----" + $nl + "4|    |" + $nl + "0|    |"
$t.Cell(5, 3).Range.Text = "88 x 64" + $nl + "  6    4" + $nl + "  ----" + $nl + "8|    |" + $nl + "8|    |"
